# Modelos de mensagens import
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the example/placeholder row (row 3): the example text now shows
# the [NOME] tag usage in both the greeting and content example cells.
$ws.Range("A3").Value = "ex.: Feliz aniversário [NOME]!"
$ws.Range("C3").Value = "ex.: É Ótimo tê-lo conosco! Feliz aniversário [NOME]!"

# Move the active selection to C3
$ws.Range("C3").Select()
